$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Insert 3 new rows before row 6, shifting existing rows (and the
# formulas on the other sheets that reference them) down by three.
$ws.Rows("6:8").Insert()

# Row 6 stays blank (it inherits the bold "header" style from row 5
# above, which already matches the target formatting).

# Row 7: new bold + yellow-highlighted section header.
$ws.Range("A7").Value = "Notes on the Brazilian adaptation"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Interior.Color = 65535
$ws.Range("B7").Interior.Color = 65535

# Row 8: explanatory note below the header (plain formatting).
$ws.Range("A8").Value = "We assumed the same currencies as the US model. "
$ws.Range("A8").Font.Bold = $false
